$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.207.89'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = '1.611.93'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9994'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.19'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3781'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3682'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '49.14'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9991'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.279'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08114'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.02%  '
$ws.Range('E13').Value = '  -3.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.634'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.647'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001274'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.41%  '
$ws.Range('D17').Value = '1.608.62'
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.59'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06806'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.608'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9979'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.09'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.25%  '
$ws.Range('D24').Value = '23.217.35'
$ws.Range('E24').Value = '  -3.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.366'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.923'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('E27').Value = '  -4.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '150.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.270'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.421'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.031'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.06%  '
$ws.Range('D33').Value = '1.787.13'
$ws.Range('E33').Value = '  -2.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9913'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07755'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02795'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.339'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.50%  '
$ws.Range('E38').Value = '  -4.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.18'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.08889'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.401'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7204'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.83%  '
$ws.Range('E43').Value = '  -4.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6638'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.317'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9974'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.980'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08023'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.10'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.170'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.82%  '
